$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.399.37"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.571.92"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'1.003"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'290.99"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.3761"
$ws.Range("D8").Value = "'49.90"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").Value = "'0.07656"
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "'21.24"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "'5.958"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "'6.908"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "1.574.15"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'0.00001135"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'90.33"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "'0.06769"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D21").Value = "'16.75"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").Value = "'6.206"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'0.5275"
$ws.Range("E23").Value = "  -4.46%  "
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "'2.452"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "22.392.53"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'2.739"
$ws.Range("E27").Value = "  -7.23%  "
$ws.Range("D28").Value = "'20.27"
$ws.Range("E28").Value = "  +2.53%  "
$ws.Range("D29").Value = "'145.33"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "'5.059"
$ws.Range("E30").Value = "  +2.88%  "
$ws.Range("D31").Value = "'125.86"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "1.747.72"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.180"
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.011"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.015"
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("D36").Value = "'10.05"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("D37").Value = "'0.08543"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").Value = "'0.02556"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.333"
$ws.Range("E40").Value = "  +6.87%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.06522"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("D43").Value = "'0.6454"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").Value = "'11.55"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "'14.07"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").Value = "'0.6026"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "'3.789"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").Value = "'1.297"
$ws.Range("E49").Value = "  +8.86%  "
$ws.Range("D50").Value = "'2.090"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").Value = "'125.16"
$ws.Range("E51").Value = "  +3.16%  "
